# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a literal text value (matches the source
    # workbook, where every data cell is stored as an inline string),
    # even though many of the strings look numeric (e.g. "1.000").
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.205.44"
Set-TextValue $ws.Range("E2") "  -0.92%  "

Set-TextValue $ws.Range("D3") "1.860.19"
Set-TextValue $ws.Range("E3") "  -0.97%  "

Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.02%  "

Set-TextValue $ws.Range("D5") "0.7067"
Set-TextValue $ws.Range("E5") "  -1.07%  "

Set-TextValue $ws.Range("D6") "240.35"
Set-TextValue $ws.Range("E6") "  -0.81%  "

Set-TextValue $ws.Range("D7") "1.001"
Set-TextValue $ws.Range("E7") "  +0.02%  "

Set-TextValue $ws.Range("D8") "0.3071"
Set-TextValue $ws.Range("E8") "  -1.49%  "

Set-TextValue $ws.Range("D9") "0.07642"
Set-TextValue $ws.Range("E9") "  -2.77%  "

Set-TextValue $ws.Range("D10") "24.71"
Set-TextValue $ws.Range("E10") "  -2.02%  "

Set-TextValue $ws.Range("D11") "0.08420"
Set-TextValue $ws.Range("E11") "  +1.83%  "

Set-TextValue $ws.Range("D12") "1.873.12"
Set-TextValue $ws.Range("E12") "  +1.51%  "

Set-TextValue $ws.Range("D13") "5.172"
Set-TextValue $ws.Range("E13") "  -2.07%  "

Set-TextValue $ws.Range("D14") "0.7083"
Set-TextValue $ws.Range("E14") "  -3.22%  "

Set-TextValue $ws.Range("D15") "91.00"
Set-TextValue $ws.Range("E15") "  -0.31%  "

Set-TextValue $ws.Range("D16") "29.212.15"
Set-TextValue $ws.Range("E16") "  -0.75%  "

Set-TextValue $ws.Range("D17") "5.922"
Set-TextValue $ws.Range("E17") "  -0.15%  "

Set-TextValue $ws.Range("D18") "243.22"
Set-TextValue $ws.Range("E18") "  -1.82%  "

Set-TextValue $ws.Range("D19") "0.000007821"
Set-TextValue $ws.Range("E19") "  -0.76%  "

Set-TextValue $ws.Range("D20") "2.113.29"
Set-TextValue $ws.Range("E20") "  +0.05%  "

Set-TextValue $ws.Range("D21") "13.10"
Set-TextValue $ws.Range("E21") "  -1.55%  "

Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("E22") "  +0.09%  "

Set-TextValue $ws.Range("D23") "7.836"
Set-TextValue $ws.Range("E23") "  -1.79%  "

Set-TextValue $ws.Range("D24") "1.001"
Set-TextValue $ws.Range("E24") "  +0.06%  "

Set-TextValue $ws.Range("D25") "0.1588"
Set-TextValue $ws.Range("E25") "  -0.33%  "

Set-TextValue $ws.Range("D26") "162.78"
Set-TextValue $ws.Range("E26") "  -0.68%  "

Set-TextValue $ws.Range("D27") "8.892"
Set-TextValue $ws.Range("E27") "  -1.41%  "

Set-TextValue $ws.Range("D28") "18.44"
Set-TextValue $ws.Range("E28") "  +0.73%  "

Set-TextValue $ws.Range("B29") "PancakeSwap"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D29") "1.499"
Set-TextValue $ws.Range("E29") "  +0.23%  "

Set-TextValue $ws.Range("B30") "Toncoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "1.316"
Set-TextValue $ws.Range("E30") "  -3.47%  "

Set-TextValue $ws.Range("D31") "4.396"
Set-TextValue $ws.Range("E31") "  +0.47%  "

Set-TextValue $ws.Range("D32") "4.221"
Set-TextValue $ws.Range("E32") "  +2.13%  "

Set-TextValue $ws.Range("D33") "0.05132"
Set-TextValue $ws.Range("E33") "  -3.52%  "

Set-TextValue $ws.Range("D34") "0.8048"
Set-TextValue $ws.Range("E34") "  +10.97%  "

Set-TextValue $ws.Range("D35") "1.912"
Set-TextValue $ws.Range("E35") "  -1.23%  "

Set-TextValue $ws.Range("D36") "1.165"
Set-TextValue $ws.Range("E36") "  -3.14%  "

Set-TextValue $ws.Range("D37") "2.682"
Set-TextValue $ws.Range("E37") "  +0.09%  "

Set-TextValue $ws.Range("D38") "0.01844"
Set-TextValue $ws.Range("E38") "  -1.39%  "

Set-TextValue $ws.Range("D39") "2.689"
Set-TextValue $ws.Range("E39") "  -1.78%  "

Set-TextValue $ws.Range("D40") "1.173.73"
Set-TextValue $ws.Range("E40") "  -7.16%  "

Set-TextValue $ws.Range("D41") "6.176"
Set-TextValue $ws.Range("E41") "  +0.74%  "

Set-TextValue $ws.Range("D42") "0.8951"
Set-TextValue $ws.Range("E42") "  -1.91%  "

Set-TextValue $ws.Range("D43") "72.77"
Set-TextValue $ws.Range("E43") "  -1.70%  "

Set-TextValue $ws.Range("D44") "1.000"
Set-TextValue $ws.Range("E44") "  -0.04%  "

Set-TextValue $ws.Range("D45") "101.72"
Set-TextValue $ws.Range("E45") "  -1.78%  "

Set-TextValue $ws.Range("D46") "2.009.08"
Set-TextValue $ws.Range("E46") "  -0.42%  "

Set-TextValue $ws.Range("D47") "0.5156"
Set-TextValue $ws.Range("E47") "  -3.38%  "

Set-TextValue $ws.Range("D48") "1.770"
Set-TextValue $ws.Range("E48") "  -0.32%  "

Set-TextValue $ws.Range("B49") "BabyDogeCoin"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D49") "0.00000000120"
Set-TextValue $ws.Range("E49") "  -0.20%  "

Set-TextValue $ws.Range("B50") "EnergySwap"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "9.251"
Set-TextValue $ws.Range("E50") "  -0.10%  "

Set-TextValue $ws.Range("B51") "Frax"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D51") "1.000"
Set-TextValue $ws.Range("E51") "  +0.38%  "
